$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.086.98'
$ws.Range("E2").Value = '  +0.51%  '
$ws.Range("D3").Value = '2.299.18'
$ws.Range("E3").Value = '  +0.30%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '300.40'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.52'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.46%  '
$ws.Range("E7").Value = '  +3.13%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("E9").Value = '  +1.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.24'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.83%  '
$ws.Range("E11").Value = '  +0.35%  '
$ws.Range("E12").Value = '  +0.70%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '17.75'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.81%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.88'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.62%  '
$ws.Range("D15").Value = '2.657.25'
$ws.Range("E15").Value = '  +0.22%  '
$ws.Range("D16").Value = '2.287.45'
$ws.Range("E16").Value = '  +0.48%  '
$ws.Range("D18").Value = '42.947.34'
$ws.Range("E19").Value = '  +4.76%  '
$ws.Range("E20").Value = '  +1.03%  '
$ws.Range("E21").Value = '  +0.38%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.18'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.81%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '237.73'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.14%  '
$ws.Range("E24").Value = '  -1.43%  '
$ws.Range("E25").Value = '  -0.43%  '
$ws.Range("E26").Value = '  -0.81%  '
$ws.Range("E27").Value = '  +0.00%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '25.05'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.53%  '
$ws.Range("E29").Value = '  +0.49%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.04'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -13.27%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '162.61'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.56%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '33.04'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.80%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.07%  '
$ws.Range("E34").Value = '  +2.65%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '18.15'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.88%  '
$ws.Range("E36").Value = '  +1.82%  '
$ws.Range("E37").Value = '  +0.36%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0694'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.08%  '
$ws.Range("E39").Value = '  +1.12%  '
$ws.Range("E40").Value = '  -0.09%  '
$ws.Range("E41").Value = '  +1.52%  '
$ws.Range("E42").Value = '  -1.98%  '
$ws.Range("D43").Value = '2.010.05'
$ws.Range("E43").Value = '  +1.96%  '
$ws.Range("E44").Value = '  -1.24%  '
$ws.Range("E45").Value = '  -6.92%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.22'
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '17.41'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.50%  '
$ws.Range("E48").Value = '  -1.25%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '54.29'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.56%  '
$ws.Range("D50").Value = '2.530.71'
$ws.Range("E50").Value = '  +0.50%  '
$ws.Range("E51").Value = '  -0.45%  '
